$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mixer2-BOM")

# Record how many of each part have now been bought/received, which
# reduces what's still needed (column E recalculates via the existing
# MAX(0, Qty-Have-Bought) formula).
$ws.Range("D10").Value = 25
$ws.Range("D21").Value = 10
$ws.Range("D22").Value = 10

# Leave the selection where work left off.
$ws.Range("D23").Select()
